$d = $word.ActiveDocument
$shapes = $d.Shapes
$out = @()
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    if ($shp.TextFrame.HasText) {
        $out += ("Shape " + $i + ": " + $shp.TextFrame.TextRange.Text)
    }
}
$out
